# Update Leave Card 1/5/2024 4:46 PM
#
# Target sheet: "2018 LEAVE CREDITS" (the bi-monthly leave-card table, Table15
# on A8:K79). The edits:
#   - A33 date shifts from 12/31/2023-ish schedule to the corrected 1/31/2024
#     end-of-period date, and so on for every following row (A34:A67), which
#     previously had no date filled in past row 34.
#   - B33 gets a new particulars entry "SP(1-0-0)" (Special Privilege Leave).
#   - K33 (REMARKS) gets the date the SP leave was filed/approved (1/5/2024),
#     formatted like the other date remarks in the column (e.g. K24:K27).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# Column A = PERIOD (end-of-month dates), one row per month from row 10 down.
$periodDates = @{
  33 = 45322
  34 = 45351
  35 = 45382
  36 = 45412
  37 = 45443
  38 = 45473
  39 = 45504
  40 = 45535
  41 = 45565
  42 = 45596
  43 = 45626
  44 = 45657
  45 = 45688
  46 = 45716
  47 = 45747
  48 = 45777
  49 = 45808
  50 = 45838
  51 = 45869
  52 = 45900
  53 = 45930
  54 = 45961
  55 = 45991
  56 = 46022
  57 = 46053
  58 = 46081
  59 = 46112
  60 = 46142
  61 = 46173
  62 = 46203
  63 = 46234
  64 = 46265
  65 = 46295
  66 = 46326
  67 = 46356
}

foreach ($row in $periodDates.Keys) {
    $ws.Cells.Item($row, 1).Value = $periodDates[$row]
}

# New PARTICULARS entry for the newly-dated row 33.
$ws.Cells.Item(33, 2).Value = "SP(1-0-0)"

# REMARKS date for row 33 (the day the leave was recorded), matching the
# date-styled remarks already used elsewhere in the column (K24:K27).
$ws.Cells.Item(33, 11).Value = 45296
$ws.Range("K24").Copy()
$ws.Cells.Item(33, 11).PasteSpecial(-4122)
$excel.CutCopyMode = $false
